# Tripadvisor New Orleans shard update:
#  1. Insert a new "State" column into the hotel_info sheet, between
#     "Hotel_Name" and "City", populated with "Louisiana" for the existing
#     data row.
#  2. Swap the tab order so "review_info" becomes the first sheet and
#     "hotel_info" becomes the second sheet.

$wb = $excel.ActiveWorkbook

# --- Step 1: add the "State" column to hotel_info -------------------------
$hotelWs = $wb.Worksheets.Item("hotel_info")

# Hotel_Name is column B, City is column C -> insert a new column at C,
# pushing City (and everything after it) one column to the right.
$hotelWs.Columns.Item(3).Insert()
$hotelWs.Cells.Item(1, 3).Value = "State"
$hotelWs.Cells.Item(2, 3).Value = "Louisiana"

# --- Step 2: reorder sheets so review_info comes first ---------------------
$reviewWs = $wb.Worksheets.Item("review_info")
$reviewWs.Move($wb.Worksheets.Item(1))
